# Update the "取得日時" (acquired datetime) timestamps in column A of the
# "ランサーズ" sheet for rows 2-10 to reflect the new append time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-11 12:42:01"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
